# Auto-generated edit script applying the Phantom_Profits market-data refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1396.7142
$ws.Range("I15").Value = 1396.7142
$ws.Range("K15").Value = 4190.142599999999
$ws.Range("M15").Value = -4021.142599999999
$ws.Range("H32").Value = 6160.25
$ws.Range("I32").Value = 2826
$ws.Range("K32").Value = 2826
$ws.Range("M32").Value = -2500
$ws.Range("H33").Value = 278.89474
$ws.Range("I33").Value = 280.5
$ws.Range("K33").Value = 280.5
$ws.Range("M33").Value = -51.5
$ws.Range("H40").Value = 1617.25
$ws.Range("I40").Value = 1529
$ws.Range("J40").Value = 1882
$ws.Range("K40").Value = 1529
$ws.Range("L40").Value = 1882
$ws.Range("M40").Value = -1354
$ws.Range("N40").Value = -2232
$ws.Range("H43").Value = 2011.625
$ws.Range("I43").Value = 2249
$ws.Range("K43").Value = 2249
$ws.Range("M43").Value = -2180
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 3765.5386
$ws.Range("I100").Value = 2995.818
$ws.Range("J100").Value = 7999
$ws.Range("K100").Value = 2995.818
$ws.Range("L100").Value = 7999
$ws.Range("M100").Value = -2454.818
$ws.Range("N100").Value = -9081
$ws.Range("H125").Value = 1601
$ws.Range("I125").Value = 1822.2
$ws.Range("J125").Value = 1478.1111
$ws.Range("K125").Value = 16399.8
$ws.Range("L125").Value = 13302.9999
$ws.Range("M125").Value = -13939.8
$ws.Range("N125").Value = -18222.9999
$ws.Range("H129").Value = 1792.2858
$ws.Range("J129").Value = 2628.375
$ws.Range("L129").Value = 7885.125
$ws.Range("N129").Value = -17885.125
$ws.Range("H131").Value = 2194.125
$ws.Range("I131").Value = 2079.1428
$ws.Range("J131").Value = 2999
$ws.Range("K131").Value = 6237.428400000001
$ws.Range("L131").Value = 8997
$ws.Range("M131").Value = -1197.428400000001
$ws.Range("N131").Value = -19077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2414.5881
$ws.Range("I32").Value = 1760.3636
$ws.Range("K32").Value = 1760.3636
$ws.Range("M32").Value = -1473.3636
$ws.Range("H74").Value = 1485.8
$ws.Range("I74").Value = 1548.2106
$ws.Range("J74").Value = 300
$ws.Range("K74").Value = 1548.2106
$ws.Range("L74").Value = 300
$ws.Range("M74").Value = -674.2106000000001
$ws.Range("N74").Value = -2048
$ws.Range("H77").Value = 1485.8
$ws.Range("I77").Value = 1548.2106
$ws.Range("K77").Value = 7741.053000000001
$ws.Range("L77").Value = 1500
$ws.Range("M77").Value = -3373.053000000001
$ws.Range("N77").Value = -10236
$ws.Range("H88").Value = 2545.4375
$ws.Range("J88").Value = 2870.1667
$ws.Range("L88").Value = 2870.1667
$ws.Range("N88").Value = -3682.1667
$ws.Range("H91").Value = 2545.4375
$ws.Range("J91").Value = 2870.1667
$ws.Range("L91").Value = 2870.1667
$ws.Range("N91").Value = -5678.1667
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 14674.75
$ws.Range("I26").Value = 14674.75
$ws.Range("K26").Value = 14674.75
$ws.Range("M26").Value = -14382.75
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H86").Value = 12038.615
$ws.Range("I86").Value = 5199.6
$ws.Range("K86").Value = 5199.6
$ws.Range("M86").Value = -4076.6
$ws.Range("H89").Value = 12038.615
$ws.Range("I89").Value = 5199.6
$ws.Range("K89").Value = 25998
$ws.Range("M89").Value = -20382
$ws.Range("H96").Value = 15374.444
$ws.Range("I96").Value = 15374.444
$ws.Range("K96").Value = 15374.444
$ws.Range("M96").Value = -12628.444
$ws.Range("H134").Value = 4088.8333
$ws.Range("I134").Value = 3922.5862
$ws.Range("K134").Value = 11767.7586
$ws.Range("M134").Value = -9232.758600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4760
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4808
$ws.Range("H58").Value = 3123.3333
$ws.Range("I58").Value = 1977.75
$ws.Range("J58").Value = 4039.8
$ws.Range("K58").Value = 1977.75
$ws.Range("L58").Value = 4039.8
$ws.Range("M58").Value = -1774.75
$ws.Range("N58").Value = -4445.8
$ws.Range("H86").Value = 3109
$ws.Range("I86").Value = 3143.125
$ws.Range("K86").Value = 3143.125
$ws.Range("M86").Value = -2020.125
$ws.Range("H89").Value = 3109
$ws.Range("I89").Value = 3143.125
$ws.Range("K89").Value = 15715.625
$ws.Range("M89").Value = -10099.625
$ws.Range("H122").Value = 9599.5
$ws.Range("I122").Value = 9772.799999999999
$ws.Range("K122").Value = 29318.4
$ws.Range("M122").Value = -26868.4
$ws.Range("H132").Value = 1514.2
$ws.Range("I132").Value = 1514.2
$ws.Range("K132").Value = 4542.6
$ws.Range("M132").Value = -2012.6
$ws.Range("H136").Value = 3123.3333
$ws.Range("I136").Value = 1977.75
$ws.Range("J136").Value = 4039.8
$ws.Range("K136").Value = 5933.25
$ws.Range("L136").Value = 12119.4
$ws.Range("M136").Value = -3383.25
$ws.Range("N136").Value = -17219.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111.4
$ws.Range("I2").Value = 113.625
$ws.Range("J2").Value = 102.5
$ws.Range("K2").Value = 681.75
$ws.Range("L2").Value = 615
$ws.Range("M2").Value = -568.75
$ws.Range("N2").Value = -841
$ws.Range("H15").Value = 307.5
$ws.Range("I15").Value = 299
$ws.Range("K15").Value = 897
$ws.Range("M15").Value = -757
$ws.Range("H17").Value = 545.6
$ws.Range("I17").Value = 289
$ws.Range("J17").Value = 609.75
$ws.Range("K17").Value = 867
$ws.Range("L17").Value = 1829.25
$ws.Range("M17").Value = -698
$ws.Range("N17").Value = -2167.25
$ws.Range("H106").Value = 13899
$ws.Range("I106").Value = 9495
$ws.Range("K106").Value = 28485
$ws.Range("M106").Value = -27539
$ws.Range("H107").Value = 113.625
$ws.Range("J107").Value = 115.166664
$ws.Range("L107").Value = 345.499992
$ws.Range("N107").Value = -4185.499992
$ws.Range("H113").Value = 4139.5835
$ws.Range("I113").Value = 3999
$ws.Range("J113").Value = 4167.7
$ws.Range("K113").Value = 11997
$ws.Range("L113").Value = 12503.1
$ws.Range("M113").Value = -9827
$ws.Range("N113").Value = -16843.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H122").Value = 3866.5789
$ws.Range("I122").Value = 3790
$ws.Range("K122").Value = 11370
$ws.Range("M122").Value = -8920
$ws.Range("H126").Value = 3623.6
$ws.Range("I126").Value = 2448
$ws.Range("J126").Value = 5387
$ws.Range("K126").Value = 7344
$ws.Range("L126").Value = 16161
$ws.Range("M126").Value = -4874
$ws.Range("N126").Value = -21101

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7060.8335
$ws.Range("I7").Value = 7033
$ws.Range("K7").Value = 7033
$ws.Range("M7").Value = -6921
$ws.Range("H40").Value = 1799.6666
$ws.Range("I40").Value = 1799.6666
$ws.Range("K40").Value = 1799.6666
$ws.Range("M40").Value = -1663.6666
$ws.Range("H46").Value = 1971.5
$ws.Range("I46").Value = 2025.8
$ws.Range("J46").Value = 1700
$ws.Range("K46").Value = 2025.8
$ws.Range("L46").Value = 1700
$ws.Range("M46").Value = -1837.8
$ws.Range("N46").Value = -2076
$ws.Range("H55").Value = 328.0909
$ws.Range("I55").Value = 237.5
$ws.Range("J55").Value = 436.8
$ws.Range("K55").Value = 237.5
$ws.Range("L55").Value = 436.8
$ws.Range("M55").Value = -64.5
$ws.Range("N55").Value = -782.8
$ws.Range("H96").Value = 31089
$ws.Range("I96").Value = 31089
$ws.Range("K96").Value = 31089
$ws.Range("M96").Value = -28343
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 7060.8335
$ws.Range("I126").Value = 7033
$ws.Range("K126").Value = 21099
$ws.Range("M126").Value = -18629
$ws.Range("H136").Value = 5993.0625
$ws.Range("I136").Value = 4449.1665
$ws.Range("K136").Value = 13347.4995
$ws.Range("M136").Value = -10797.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2408.875
$ws.Range("I132").Value = 2198.946
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 6596.838
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -4066.838
$ws.Range("N132").Value = -20054

Write-Host "Applied all Phantom_Profits cell updates"
